# RDCC-5182 SRD file version validation changes
# Adds a new "VERSION" worksheet as the first sheet in the workbook,
# containing the file version info in cells A6/B6, and makes it the
# active sheet.

$wb = $excel.ActiveWorkbook

# Add a new worksheet before the first existing sheet ("Staff Data")
$firstSheet = $wb.Worksheets.Item(1)
$versionSheet = $wb.Worksheets.Add($firstSheet)
$versionSheet.Name = "VERSION"

# Populate the version info cells
$versionSheet.Range("A6").Value = "File version"
$versionSheet.Range("B6").Value = "vx.xx"

# Make the new sheet the active one, with B6 selected
$versionSheet.Activate()
$versionSheet.Range("B6").Select()
